$d = $word.ActiveDocument

# The "Role"/"Responsibilities" table is the first table in the document.
$t = $d.Tables.Item(1)

# Add a new row at the end of the table for the "Developers" role.
$newRow = $t.Rows.Add()

$cell1 = $newRow.Cells.Item(1)
$cell1.Range.Text = "Developers"
$cell1.Range.Font.Bold = 1

$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "Junior developers building the LMS"
